$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data, pushing all student rows
# down by one (the data itself is unchanged, just relocated).
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Student ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Year"
$ws.Range("D1").Value = "Section"
$ws.Range("E1").Value = "Contact"
$ws.Range("F1").Value = "Number of Borrowed Books"
